$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "322.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.79%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.88%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.884"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "11.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08028"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.03%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.651"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.18%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.936"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.28%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9324"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.28%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1234"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.33%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1963"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.25%"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.745"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "21.03%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09205"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.53%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03545"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.17%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09566"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.19%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001292"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-7.24%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006268"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.08%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.347"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.39%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.572"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.98%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.952"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.22%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.00%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1418"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "7.98%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.15%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.66%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001261"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.13%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004391"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.80%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-11.63%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.03%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02417"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.86%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05226"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.24%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007450"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.23%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009374"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.62%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1405"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.90%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002119"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.86%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01122"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "37.73%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006736"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.07%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003001"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-7.40%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.05%"
